# Fixes total coal values in historical extension
#
# The "Power" sector lookup values used a shared string with stray
# leading/trailing whitespace (" Power       "), which made it fail to
# match against other, clean "Power" strings elsewhere in the coal
# aggregation pipeline. Replace the affected cells with a clean "Power"
# value (Excel will intern it as a new shared string).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Bond_to_ext
$ws2 = $wb.Worksheets.Item(2)   # CEDS_to_ext

# --- Data fix: replace the whitespace-padded " Power " lookups with a
#     clean "Power" string -----------------------------------------------
$ws1.Range("C23").Value = "Power"

$ws2.Range("B2").Value = "Power"
$ws2.Range("B3").Value = "Power"
$ws2.Range("B4").Value = "Power"

# --- View state left behind by the author's last editing session --------

# Bond_to_ext: scrolled/selected around C23 (no longer the active tab)
$ws1.Activate()
$ws1.Range("C23").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1

# CEDS_to_ext: becomes the active tab, selection parked at A6
$ws2.Activate()
$ws2.Range("A6").Select()
